{"js": "// Update the date line and every \"a\u00f7b=c, d\" answer cell in the worksheet\n// table to the new values, per the commit's regenerated output.\nconst replacements = [\n  [\"2024-12-22 Sunday\", \"2024-12-23 Monday\"],\n  [\"17\u00f79=1, 8\", \"42\u00f75=8, 2\"],\n  [\"10\u00f79=1, 1\", \"83\u00f78=10, 3\"],\n  [\"18\u00f73=6, 0\", \"20\u00f76=3, 2\"],\n  [\"22\u00f79=2, 4\", \"22\u00f76=3, 4\"],\n  [\"16\u00f73=5, 1\", \"85\u00f79=9, 4\"],\n  [\"29\u00f73=9, 2\", \"77\u00f76=12, 5\"],\n  [\"96\u00f77=13, 5\", \"39\u00f79=4, 3\"],\n  [\"86\u00f72=43, 0\", \"88\u00f74=22, 0\"],\n  [\"55\u00f77=7, 6\", \"50\u00f74=12, 2\"],\n  [\"61\u00f76=10, 1\", \"59\u00f74=14, 3\"],\n  [\"45\u00f79=5, 0\", \"51\u00f76=8, 3\"],\n  [\"36\u00f78=4, 4\", \"51\u00f75=10, 1\"],\n  [\"47\u00f76=7, 5\", \"10\u00f72=5, 0\"],\n  [\"60\u00f72=30, 0\", \"80\u00f78=10, 0\"],\n  [\"81\u00f76=13, 3\", \"56\u00f79=6, 2\"],\n  [\"88\u00f75=17, 3\", \"64\u00f79=7, 1\"],\n  [\"62\u00f74=15, 2\", \"33\u00f75=6, 3\"],\n  [\"80\u00f76=13, 2\", \"81\u00f75=16, 1\"],\n  [\"36\u00f73=12, 0\", \"21\u00f72=10, 1\"],\n  [\"48\u00f79=5, 3\", \"70\u00f75=14, 0\"],\n  [\"41\u00f79=4, 5\", \"84\u00f75=16, 4\"],\n  [\"24\u00f77=3, 3\", \"48\u00f77=6, 6\"],\n  [\"21\u00f75=4, 1\", \"47\u00f77=6, 5\"],\n  [\"87\u00f77=12, 3\", \"18\u00f79=2, 0\"],\n  [\"48\u00f73=16, 0\", \"65\u00f74=16, 1\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, {\n    matchCase: true,\n    matchWholeWord: false,\n  });\n  results.load(\"text\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date line and every \"a\u00f7b=c, d\" answer cell in the worksheet\n# table to the new values, per the commit's regenerated output.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2024-12-22 Sunday\", \"2024-12-23 Monday\"),\n    @(\"17\u00f79=1, 8\", \"42\u00f75=8, 2\"),\n    @(\"10\u00f79=1, 1\", \"83\u00f78=10, 3\"),\n    @(\"18\u00f73=6, 0\", \"20\u00f76=3, 2\"),\n    @(\"22\u00f79=2, 4\", \"22\u00f76=3, 4\"),\n    @(\"16\u00f73=5, 1\", \"85\u00f79=9, 4\"),\n    @(\"29\u00f73=9, 2\", \"77\u00f76=12, 5\"),\n    @(\"96\u00f77=13, 5\", \"39\u00f79=4, 3\"),\n    @(\"86\u00f72=43, 0\", \"88\u00f74=22, 0\"),\n    @(\"55\u00f77=7, 6\", \"50\u00f74=12, 2\"),\n    @(\"61\u00f76=10, 1\", \"59\u00f74=14, 3\"),\n    @(\"45\u00f79=5, 0\", \"51\u00f76=8, 3\"),\n    @(\"36\u00f78=4, 4\", \"51\u00f75=10, 1\"),\n    @(\"47\u00f76=7, 5\", \"10\u00f72=5, 0\"),\n    @(\"60\u00f72=30, 0\", \"80\u00f78=10, 0\"),\n    @(\"81\u00f76=13, 3\", \"56\u00f79=6, 2\"),\n    @(\"88\u00f75=17, 3\", \"64\u00f79=7, 1\"),\n    @(\"62\u00f74=15, 2\", \"33\u00f75=6, 3\"),\n    @(\"80\u00f76=13, 2\", \"81\u00f75=16, 1\"),\n    @(\"36\u00f73=12, 0\", \"21\u00f72=10, 1\"),\n    @(\"48\u00f79=5, 3\", \"70\u00f75=14, 0\"),\n    @(\"41\u00f79=4, 5\", \"84\u00f75=16, 4\"),\n    @(\"24\u00f77=3, 3\", \"48\u00f77=6, 6\"),\n    @(\"21\u00f75=4, 1\", \"47\u00f77=6, 5\"),\n    @(\"87\u00f77=12, 3\", \"18\u00f79=2, 0\"),\n    @(\"48\u00f73=16, 0\", \"65\u00f74=16, 1\")\n)\n\nforeach ($pair in $pairs) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $found = $find.Execute($pair[0], $false, $false, $false, $false, $false, $true, 1, $false, $pair[1], 2)\n    if (-not $found) {\n        throw \"Text not found: $($pair[0])\"\n    }\n}\n"}
